$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "3 - Time Greeting" -> "3 - Ask Services"
$ws.Range("A3").Value = "3 - Ask Services"

# New column B notes for rows 1 and 2
$ws.Range("B1").Value = "Used Hi for all injection words (Hey, Hi,..)"
$ws.Range("B2").Value = "Replaced usergreetings -> Good Morning"

# Rows 4 and 5 no longer hold the old "4 - ..." / "5 - ..." labels
$ws.Range("A4").ClearContents()
$ws.Range("A5").ClearContents()

# Widen column B to match the new note column
$ws.Columns.Item(2).ColumnWidth = 41.7

# Move the active selection to B4
$ws.Range("B4").Select() | Out-Null
